# Add hw9 / hw10 grade columns (J, K) to the homework grades sheet and
# backfill grades for hw7 (H) / hw8 (I) / hw9 (J) for all students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -------------------------------------------------------
$ws.Range("J1").Value = "hw9"
$ws.Range("K1").Value = "hw10"

# ---- Grade values (hw7=H, hw8=I, hw9=J) --------------------------------
# row -> H, I, J
$grades = @{
    2  = @(99, 99, 99)
    3  = @(99, 97, 98)
    4  = @(98, 99, 99)
    5  = @(93, 98, 94)
    6  = @(94, 0, 97)
    7  = @(100, 97, 99)
    8  = @(92, 100, 89)
    9  = @(100, 99, 99)
    10 = @(100, 100, 100)
    11 = @(100, 99, 100)
    12 = @(98, 99, 98)
    13 = @(92, 89, 90)
    14 = @(99, 99, 98)
    15 = @(98, 95, 99)
    16 = @(100, 99, 100)
}

foreach ($row in $grades.Keys) {
    $vals = $grades[$row]
    $ws.Range("H$row").Value = $vals[0]
    $ws.Range("I$row").Value = $vals[1]
    $ws.Range("J$row").Value = $vals[2]
}

# ---- Formatting ---------------------------------------------------------
# Most of the new/touched grade cells pick up the plain bordered style
# already used elsewhere in the sheet (same as cell B4). The lone
# "did not turn in" (0) grade, I6, picks up the red "did not turn in"
# style already used for B6.
$ws.Range("B4").Copy()
$ws.Range("H2:J5").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("H7:J16").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Column widths for the two new columns -----------------------------
$ws.Columns("J").ColumnWidth = 3.83
$ws.Columns("K").ColumnWidth = 4.83

# ---- Selection matches the saved state in the source workbook ----------
$ws.Range("I17").Select()
